# Repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value  = -8
$ws.Range("F9").Value  = 1
$ws.Range("F13").Value = -3
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = -1
